$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the value in B2 (bl) from 20 to 10, per new GAS hierarchy
$ws.Range("B2").Value = 10

# Update the active cell selection to B3
$ws.Range("B3").Select()
